$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.962.25"
$ws.Range("E2").Value = "  +3.00%  "
$ws.Range("D3").Value = "3.802.43"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "699.84"
$ws.Range("E5").Value = "  +8.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.62"
$ws.Range("E6").Value = "  +4.38%  "
$ws.Range("D7").Value = "3.801.92"
$ws.Range("E7").Value = "  +1.02%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("E10").Value = "  +3.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.40"
$ws.Range("E11").Value = "  +6.64%  "
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000258"
$ws.Range("E13").Value = "  +8.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.44"
$ws.Range("E14").Value = "  +4.62%  "
$ws.Range("D15").Value = "4.442.79"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "3.793.59"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").Value = "70.926.67"
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.90"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("E19").Value = "  +3.12%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.07"
$ws.Range("E21").Value = "  +15.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "482.97"
$ws.Range("E22").Value = "  +2.20%  "
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.44"
$ws.Range("E24").Value = "  +3.21%  "
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  +2.69%  "
$ws.Range("E27").Value = "  +4.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.49"
$ws.Range("E28").Value = "  +4.24%  "
$ws.Range("D29").Value = "3.952.20"
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  +16.19%  "
$ws.Range("E32").Value = "  +6.28%  "
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("E34").Value = "  +6.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "29.50"
$ws.Range("E35").Value = "  +3.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.24"
$ws.Range("E36").Value = "  +4.39%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +2.53%  "
$ws.Range("E39").Value = "  +7.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.04"
$ws.Range("E40").Value = "  +4.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.21"
$ws.Range("E41").Value = "  +12.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.000328"
$ws.Range("E42").Value = "  +23.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.976"
$ws.Range("E43").Value = "  +1.97%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "162.36"
$ws.Range("E46").Value = "  +4.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.24"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.84"
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("E49").Value = "  +3.04%  "
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("E51").Value = "  +2.78%  "
